# Updates the crypto price table with refreshed values from coinranking.com
# (rows for some coins were also re-ordered/re-paired with new tokens).
# NumberFormat '@' (Text) is applied before writing numeric-looking price
# strings so Excel stores/keeps the exact literal text (e.g. "3.560",
# "0.0001500") instead of silently re-parsing it into a Double and losing
# trailing zeros, matching the original inline-string cell content.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '266.65'
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '21.34'
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '6.108'
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '0.06108'
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '3.560'
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '6.489'
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '1.357'
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.8205'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.01339'
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.1596'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.08031'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.03456'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.03214'
$ws.Range('B15').Value = 'BitMartToken'
$ws.Range('C15').Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.09219'
$ws.Range('E15').Value = '14BitMartTokenBMX'
$ws.Range('B16').Value = 'MCDex'
$ws.Range('C16').Value = 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '3.732'
$ws.Range('E16').Value = '15MCDexMCB'
$ws.Range('B17').Value = 'BitForexToken'
$ws.Range('C17').Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.001620'
$ws.Range('E17').Value = '16BitForexTokenBF'
$ws.Range('B18').Value = 'CoinExToken'
$ws.Range('C18').Value = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.04631'
$ws.Range('E18').Value = '17CoinExTokenCET'
$ws.Range('B19').Value = 'TigerCash'
$ws.Range('C19').Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.006331'
$ws.Range('E19').Value = '18TigerCashTCH'
$ws.Range('B20').Value = 'HotbitToken'
$ws.Range('C20').Value = 'https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.006149'
$ws.Range('E20').Value = '19HotbitTokenHTB'
$ws.Range('B21').Value = 'BitKan'
$ws.Range('C21').Value = 'https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.001068'
$ws.Range('E21').Value = '20BitKanKAN'
$ws.Range('B22').Value = 'NitroEx'
$ws.Range('C22').Value = 'https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.0001500'
$ws.Range('E22').Value = '21NitroExNTX'
$ws.Range('B23').Value = 'LEO'
$ws.Range('C23').Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '3.728'
$ws.Range('E23').Value = '22LEOLEO'
$ws.Range('B24').Value = 'BTSEToken'
$ws.Range('C24').Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.270'
$ws.Range('E24').Value = '23BTSETokenBTSE'
$ws.Range('B25').Value = 'BitpandaEcosystemToken'
$ws.Range('C25').Value = 'https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.3314'
$ws.Range('E25').Value = '24BitpandaEcosystemTokenBEST'
$ws.Range('B26').Value = 'ProBitToken'
$ws.Range('C26').Value = 'https://coinranking.com/coin/lQP4d6T2+probittoken-prob'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.1242'
$ws.Range('E26').Value = '25ProBitTokenPROB'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.04600'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.006963'
$ws.Range('B42').Value = 'CEJI'
$ws.Range('C42').Value = 'https://coinranking.com/coin/SbKjCVJCh+ceji-ceji'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.003765'
$ws.Range('E42').Value = '41CEJICEJI'
$ws.Range('B43').Value = 'BKEXToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.1116'
$ws.Range('E43').Value = '42BKEXTokenBKK'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.01058'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.00005912'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.00001900'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.01240'
